# "Add files via upload" - re-upload of the "Prueba escritorio" workbook with
# the log matrix on Hoja1 trimmed (some cells flipped from 1 -> 0) and the
# last selected cell on that sheet moved from A15 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the log/result grid (rows 6, 7, 9, 10) ---
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

$ws.Range("H10").Value = 0

# --- Restore/move the active selection on Hoja1 to C5 ---
$ws.Activate()
[void]$ws.Range("C5").Select()
